# Updated cryptos list values (price + 1h volume change) per source diff,
# refreshed by the scheduled GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of these cells are stored as literal text in the workbook (prices such as
# "27.694.32" / "1.886.93" use dots as thousands separators and are not valid
# numbers, and the Volume(1h) column is a padded "  +0.44%  " string). Force every
# touched cell to Text format *before* assigning so Excel does not opportunistically
# reinterpret numeric-looking strings (e.g. "5.470", "1.200") as numbers and silently
# drop significant trailing zeros.
$textCells = @(
'D2', 'E2', 'D3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'E16', 'D17', 'E17', 'D18', 'E18', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'D29', 'E29', 'D30', 'E30', 'D31', 'E31', 'D32', 'D33', 'E33', 'D34', 'E34', 'D35', 'E35', 'D36', 'E36', 'D37', 'E37', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'D50', 'E50', 'B51', 'C51', 'D51', 'E51'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.701.48'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '1.875.60'
$ws.Range('D4').Value = '1.014'
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '336.44'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('D6').Value = '1.013'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').Value = '0.4663'
$ws.Range('E7').Value = '  -1.44%  '
$ws.Range('D8').Value = '0.3941'
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('D9').Value = '45.92'
$ws.Range('E9').Value = '  -3.94%  '
$ws.Range('D10').Value = '0.07995'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('E11').Value = '  -1.87%  '
$ws.Range('D12').Value = '21.82'
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('D13').Value = '1.888.85'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').Value = '5.979'
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').Value = '7.254'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = '88.96'
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('D18').Value = '0.06748'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('D20').Value = '17.29'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '1.011'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = '27.697.42'
$ws.Range('E22').Value = '  -0.55%  '
$ws.Range('D23').Value = '5.470'
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('D24').Value = '10.96'
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').Value = '2.308'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').Value = '2.106.96'
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('D27').Value = '159.38'
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').Value = '19.77'
$ws.Range('D29').Value = '2.153'
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('D30').Value = '5.482'
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('D31').Value = '121.53'
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('D32').Value = '0.9790'
$ws.Range('D33').Value = '0.09456'
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('D34').Value = '3.633'
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('D35').Value = '5.326'
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('D36').Value = '1.350'
$ws.Range('E36').Value = '  -6.94%  '
$ws.Range('D37').Value = '0.06063'
$ws.Range('E37').Value = '  -1.55%  '
$ws.Range('D38').Value = '0.02240'
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').Value = '8.338'
$ws.Range('E39').Value = '  +3.52%  '
$ws.Range('D40').Value = '1.200'
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('D41').Value = '1.013'
$ws.Range('E41').Value = '  +0.50%  '
$ws.Range('D42').Value = '0.5963'
$ws.Range('D43').Value = '0.1882'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').Value = '10.37'
$ws.Range('E44').Value = '  +0.71%  '
$ws.Range('D45').Value = '1.251'
$ws.Range('E45').Value = '  -0.85%  '
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('D47').Value = '12.24'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').Value = '1.935'
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D49').Value = '0.06763'
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('D50').Value = '111.96'
$ws.Range('E50').Value = '  -1.34%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = '1.060'
$ws.Range('E51').Value = '  -1.08%  '
